$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Update the form_id setting's value from "refrigerators" to "refrigerators_init"
# (the table_id setting keeps its "refrigerators" value)
$ws.Cells.Item(2, 2).Value = "refrigerators_init"
